$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (D): values look like plain numbers (e.g. "568.50") but
# must stay literal text (matching the source data, which stores them as
# strings, sometimes using "." as a thousands separator, e.g. "65.408.96").
# Forcing the cell to Text format before assignment keeps Excel from
# reinterpreting the string as a number.
$priceUpdates = @{
    'D2' = '65.408.96'
    'D3' = '2.936.20'
    'D5' = '568.50'
    'D6' = '158.42'
    'D8' = '0.515'
    'D9' = '2.932.87'
    'D10' = '6.74'
    'D12' = '0.459'
    'D14' = '34.32'
    'D16' = '65.366.19'
    'D17' = '3.420.89'
    'D18' = '6.98'
    'D19' = '2.931.09'
    'D20' = '15.68'
    'D21' = '444.03'
    'D22' = '0.691'
    'D23' = '7.27'
    'D24' = '82.21'
    'D25' = '2.24'
    'D26' = '12.11'
    'D27' = '10.05'
    'D29' = '8.03'
    'D31' = '2.58'
    'D32' = '0.0000100'
    'D36' = '0.973'
    'D37' = '5.74'
    'D38' = '49.63'
    'D39' = '44.90'
    'D40' = '1.98'
    'D44' = '8.48'
    'D45' = '383.13'
    'D46' = '0.0352'
    'D47' = '2.697.63'
    'D48' = '133.46'
    'D51' = '0.107'
}

foreach ($cell in $priceUpdates.Keys) {
    $range = $ws.Range($cell)
    $range.NumberFormat = "@"
    $range.Value = $priceUpdates[$cell]
}

# --- Volume(1h) column (E): padded percentage strings (e.g. "  -1.15%  ").
# These are not valid numeric literals as-is, so a plain assignment keeps them
# as text without needing to touch the cell number format.
$volumeUpdates = @{
    'E2' = '  -1.15%  '
    'E3' = '  -2.72%  '
    'E4' = '  -0.05%  '
    'E5' = '  -2.84%  '
    'E6' = '  +1.45%  '
    'E7' = '  -0.01%  '
    'E8' = '  -0.60%  '
    'E9' = '  -2.69%  '
    'E10' = '  -3.56%  '
    'E11' = '  -3.58%  '
    'E12' = '  +1.47%  '
    'E13' = '  -2.89%  '
    'E14' = '  -0.96%  '
    'E15' = '  -0.83%  '
    'E16' = '  -1.23%  '
    'E17' = '  -2.75%  '
    'E18' = '  +0.30%  '
    'E19' = '  -3.02%  '
    'E20' = '  +12.98%  '
    'E21' = '  -3.99%  '
    'E22' = '  +0.80%  '
    'E23' = '  -1.33%  '
    'E24' = '  +0.11%  '
    'E25' = '  -1.39%  '
    'E26' = '  -3.05%  '
    'E27' = '  -6.19%  '
    'E28' = '  +0.10%  '
    'E29' = '  -0.12%  '
    'E30' = '  -0.07%  '
    'E31' = '  -1.49%  '
    'E32' = '  -4.10%  '
    'E33' = '  +0.03%  '
    'E34' = '  -0.20%  '
    'E35' = '  +0.00%  '
    'E36' = '  -2.34%  '
    'E37' = '  -1.35%  '
    'E38' = '  +0.22%  '
    'E39' = '  +2.82%  '
    'E40' = '  -8.84%  '
    'E41' = '  -1.02%  '
    'E42' = '  -2.10%  '
    'E43' = '  -7.99%  '
    'E44' = '  +0.42%  '
    'E45' = '  -2.77%  '
    'E46' = '  -0.78%  '
    'E47' = '  -3.69%  '
    'E48' = '  -0.34%  '
    'E50' = '  +4.71%  '
    'E51' = '  +0.18%  '
}

foreach ($cell in $volumeUpdates.Keys) {
    $ws.Range($cell).Value = $volumeUpdates[$cell]
}
